# Applies the table style change on slide 5's table:
#   tableStyleId {860E92E5-11E5-48CB-952C-1525C6F75D2B} -> {4961DFB1-00C7-43F5-8F93-9A1F6392C039}
#
# Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") has 3 top-level shapes; the 2nd
# one is the p:graphicFrame hosting the a:tbl (the only table in the deck that
# references the old style id). PowerPoint's object model only allows
# re-styling a table through Table.ApplyStyle("{GUID}") - a direct assignment
# to Table.Style throws ("Table styles cannot be assigned through a property -
# call Table.ApplyStyle(\"{GUID}\") instead"), which is what we use below.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{4961DFB1-00C7-43F5-8F93-9A1F6392C039}")
